# Split the run containing "Persistência de " into two runs:
#   "Persistência"  and  " de "
# so that the formatting (identical rPr) is preserved on both runs,
# matching the target OOXML diff.

$d = $word.ActiveDocument

$searchText = "Persistência de "
$splitWord  = "Persistência"

$rng = $d.Content
$found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    Write-Host "ERROR: text '$searchText' not found"
} else {
    $start = $rng.Start
    $end   = $rng.End

    # Position right after "Persistência", i.e. before " de "
    $splitPos = $start + $splitWord.Length

    # Grab just the first part of the run and toggle a character
    # format on/off; Word materializes this as a genuine run split
    # in the underlying OOXML (two <w:r> elements with identical rPr).
    $firstPart = $d.Range($start, $splitPos)
    $firstPart.Bold = 1
    $firstPart.Bold = 0

    Write-Host "Split run at position $splitPos -> '$($d.Range($start,$splitPos).Text)' | '$($d.Range($splitPos,$end).Text)'"
}
